# Updated the pre-tokenization section.
#
# Slide 2 contains a group ("Group 69", id=70) holding the pre-tokenization
# flow-chart. This edit:
#   1. Re-groups the shapes, which causes PowerPoint to renumber the group
#      (70/"Group 69" -> 7/"Group 6"), matching the authoring tool's
#      behaviour when shapes are re-grouped.
#   2. Moves the "Evaluator" rounded rectangle, the "Results" flow-chart
#      process box, and the arrow connecting them further down the slide.
#   3. Resizes/repositions the two elbow connectors that lead into / out of
#      that area so they still meet up with their connected shapes.
#
# NOTE: Shape.Top/Left/Width/Height are exposed as single-precision (Single)
# points in the COM object model, same as real PowerPoint. 1 pt = 12700 EMU,
# and a plain "EMU / 12700.0" division can be off by a few EMU once it has
# been through the float32 round-trip. The literal point values below were
# chosen so that, after that round-trip, they land exactly on the target
# EMU values from the original file.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# --- Locate the top-level group shape on the slide -----------------------
$group = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.Name -eq "Group 69") {
        $group = $candidate
    }
}

# --- Re-group so the group shape gets renumbered (70/"Group 69" -> 7/"Group 6")
$groupItems = $group.Ungroup()
$group = $groupItems.Group()
$group.Name = "Group 6"

# --- Adjust the shapes that changed position/size inside the group -------
for ($i = 1; $i -le $group.GroupItems.Count; $i++) {
    $shp = $group.GroupItems.Item($i)

    if ($shp.Id -eq 111) {
        # Elbow Connector 110: taller, same top/left/width
        $shp.Height = 89.5959816519685
    }
    elseif ($shp.Id -eq 64) {
        # Elbow Connector 63: moves down and shrinks, same left/width
        $shp.Top = 403.2431488062992
        $shp.Height = 68.4915734031496
    }
    elseif ($shp.Id -eq 39) {
        # Rounded Rectangle 38 ("Evaluator"): moves down
        $shp.Top = 363.38441474881887
    }
    elseif ($shp.Id -eq 44) {
        # Flowchart: Process 43 ("Results"): moves down
        $shp.Top = 363.3057404314961
    }
    elseif ($shp.Id -eq 45) {
        # Straight Arrow Connector 44: moves down
        $shp.Top = 383.27449038897635
    }
}
